$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (30 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 201
$ws.Range("I9").Value = 200
$ws.Range("J9").Value = 202
$ws.Range("K9").Value = 200
$ws.Range("L9").Value = 202
$ws.Range("M9").Value = -31
$ws.Range("N9").Value = -540
$ws.Range("H32").Value = 3953.3635
$ws.Range("I32").Value = 3873
$ws.Range("K32").Value = 3873
$ws.Range("M32").Value = -3547
$ws.Range("H62").Value = 7743.467
$ws.Range("I62").Value = 7396.4614
$ws.Range("K62").Value = 7396.4614
$ws.Range("M62").Value = -6772.4614
$ws.Range("H65").Value = 7743.467
$ws.Range("I65").Value = 7396.4614
$ws.Range("K65").Value = 36982.307
$ws.Range("M65").Value = -33862.307
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H138").Value = 2702.4
$ws.Range("I138").Value = 2619.9412
$ws.Range("J138").Value = 2763.348
$ws.Range("K138").Value = 7859.823600000001
$ws.Range("L138").Value = 8290.044
$ws.Range("M138").Value = -2719.823600000001
$ws.Range("N138").Value = -18570.044

# --- Sheet: ARM (22 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1369.375
$ws.Range("I2").Value = 1369.375
$ws.Range("K2").Value = 1369.375
$ws.Range("M2").Value = -1256.375
$ws.Range("H32").Value = 21321.6
$ws.Range("I32").Value = 24842.582
$ws.Range("J32").Value = 8704.75
$ws.Range("K32").Value = 24842.582
$ws.Range("L32").Value = 8704.75
$ws.Range("M32").Value = -24555.582
$ws.Range("N32").Value = -9278.75
$ws.Range("H116").Value = 1369.375
$ws.Range("I116").Value = 1369.375
$ws.Range("K116").Value = 1369.375
$ws.Range("M116").Value = 924.625
$ws.Range("H122").Value = 1346.2
$ws.Range("I122").Value = 1356.3823
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4069.1469
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -1619.1469
$ws.Range("N122").Value = -7900

# --- Sheet: BSM (20 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1369.375
$ws.Range("I3").Value = 1369.375
$ws.Range("K3").Value = 1369.375
$ws.Range("M3").Value = -1255.375
$ws.Range("H22").Value = 56538.777
$ws.Range("I22").Value = 67644.47
$ws.Range("K22").Value = 67644.47
$ws.Range("M22").Value = -67471.47
$ws.Range("H100").Value = 28339.25
$ws.Range("J100").Value = 28339.25
$ws.Range("L100").Value = 28339.25
$ws.Range("N100").Value = -30503.25
$ws.Range("H105").Value = 120697.94
$ws.Range("I105").Value = 2919.6428
$ws.Range("K105").Value = 2919.6428
$ws.Range("M105").Value = -1172.6428
$ws.Range("H112").Value = 149695
$ws.Range("J112").Value = 149695
$ws.Range("L112").Value = 149695
$ws.Range("N112").Value = -152649

# --- Sheet: CRP (18 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1444
$ws.Range("I22").Value = 1069.6
$ws.Range("J22").Value = 1614.1818
$ws.Range("K22").Value = 1069.6
$ws.Range("L22").Value = 1614.1818
$ws.Range("M22").Value = -719.5999999999999
$ws.Range("N22").Value = -2314.1818
$ws.Range("H105").Value = 2433.1428
$ws.Range("I105").Value = 2236.8333
$ws.Range("J105").Value = 3611
$ws.Range("K105").Value = 2236.8333
$ws.Range("L105").Value = 3611
$ws.Range("M105").Value = -489.8332999999998
$ws.Range("N105").Value = -7105
$ws.Range("H122").Value = 3467
$ws.Range("I122").Value = 3090.8572
$ws.Range("K122").Value = 9272.571599999999
$ws.Range("M122").Value = -6822.571599999999

# --- Sheet: CUL (7 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2374.25
$ws.Range("I51").Value = 2374.25
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 7122.75
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -6662.75
$ws.Range("N51").ClearContents()

# --- Sheet: GSM (37 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 59998.5
$ws.Range("J63").Value = 59998.5
$ws.Range("L63").Value = 59998.5
$ws.Range("N63").Value = -61370.5
$ws.Range("H66").Value = 59998.5
$ws.Range("J66").Value = 59998.5
$ws.Range("L66").Value = 179995.5
$ws.Range("N66").Value = -186859.5
$ws.Range("H80").Value = 2858.6155
$ws.Range("I80").Value = 2681.2
$ws.Range("J80").Value = 3450
$ws.Range("K80").Value = 2681.2
$ws.Range("L80").Value = 3450
$ws.Range("M80").Value = -1683.2
$ws.Range("N80").Value = -5446
$ws.Range("H83").Value = 2858.6155
$ws.Range("I83").Value = 2681.2
$ws.Range("J83").Value = 3450
$ws.Range("K83").Value = 13406
$ws.Range("L83").Value = 17250
$ws.Range("M83").Value = -8414
$ws.Range("N83").Value = -27234
$ws.Range("H102").Value = 4157.778
$ws.Range("J102").Value = 3191.875
$ws.Range("L102").Value = 3191.875
$ws.Range("N102").Value = -6435.875
$ws.Range("H113").Value = 171882.92
$ws.Range("I113").Value = 170599.83
$ws.Range("J113").Value = 173166
$ws.Range("K113").Value = 170599.83
$ws.Range("L113").Value = 173166
$ws.Range("M113").Value = -168429.83
$ws.Range("N113").Value = -177506
$ws.Range("H126").Value = 7091.364
$ws.Range("I126").Value = 6742.8667
$ws.Range("K126").Value = 20228.6001
$ws.Range("M126").Value = -17758.6001

# --- Sheet: LTW (29 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2433.9285
$ws.Range("I7").Value = 2313.4614
$ws.Range("K7").Value = 2313.4614
$ws.Range("M7").Value = -2201.4614
$ws.Range("H22").Value = 86945.46000000001
$ws.Range("I22").Value = 278074.75
$ws.Range("J22").Value = 1999.1111
$ws.Range("K22").Value = 278074.75
$ws.Range("L22").Value = 1999.1111
$ws.Range("M22").Value = -277779.75
$ws.Range("N22").Value = -2589.1111
$ws.Range("H27").Value = 86945.46000000001
$ws.Range("I27").Value = 278074.75
$ws.Range("J27").Value = 1999.1111
$ws.Range("K27").Value = 278074.75
$ws.Range("L27").Value = 1999.1111
$ws.Range("M27").Value = -277967.75
$ws.Range("N27").Value = -2213.1111
$ws.Range("H93").Value = 3748.2
$ws.Range("I93").Value = 2498.25
$ws.Range("J93").Value = 4581.5
$ws.Range("K93").Value = 2498.25
$ws.Range("L93").Value = 4581.5
$ws.Range("M93").Value = -1250.25
$ws.Range("N93").Value = -7077.5
$ws.Range("H126").Value = 2433.9285
$ws.Range("I126").Value = 2313.4614
$ws.Range("K126").Value = 6940.3842
$ws.Range("M126").Value = -4470.3842

# --- Sheet: WVR (34 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9999
$ws.Range("I15").Value = 9999
$ws.Range("K15").Value = 9999
$ws.Range("M15").Value = -9711
$ws.Range("H54").Value = 7073.9
$ws.Range("I54").Value = 3848.125
$ws.Range("J54").Value = 19977
$ws.Range("K54").Value = 3848.125
$ws.Range("L54").Value = 19977
$ws.Range("M54").Value = -3328.125
$ws.Range("N54").Value = -21017
$ws.Range("H63").Value = 23748.5
$ws.Range("J63").Value = 23748.5
$ws.Range("L63").Value = 23748.5
$ws.Range("N63").Value = -24996.5
$ws.Range("H66").Value = 23748.5
$ws.Range("J66").Value = 23748.5
$ws.Range("L66").Value = 71245.5
$ws.Range("N66").Value = -77485.5
$ws.Range("H97").Value = 37928.75
$ws.Range("J97").Value = 37928.75
$ws.Range("L97").Value = 37928.75
$ws.Range("N97").Value = -39910.75
$ws.Range("H102").Value = 69998.5
$ws.Range("J102").Value = 69998.5
$ws.Range("L102").Value = 69998.5
$ws.Range("N102").Value = -76488.5
$ws.Range("H122").Value = 2587.5789
$ws.Range("I122").Value = 2234.5334
$ws.Range("J122").Value = 3911.5
$ws.Range("K122").Value = 6703.600199999999
$ws.Range("L122").Value = 11734.5
$ws.Range("M122").Value = -4253.600199999999
$ws.Range("N122").Value = -16634.5
